$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 58.53846
$ws.Range("I8").Value = 58.53846
$ws.Range("K8").Value = 175.61538
$ws.Range("M8").Value = -36.61538000000002
$ws.Range("H11").Value = 124.5
$ws.Range("I11").Value = 124.5
$ws.Range("K11").Value = 124.5
$ws.Range("M11").Value = 15.5
$ws.Range("H62").Value = 12680.704
$ws.Range("J62").Value = 5995
$ws.Range("L62").Value = 5995
$ws.Range("N62").Value = -7243
$ws.Range("H64").Value = 65125
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 168666.67
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 168666.67
$ws.Range("M64").Value = -2752
$ws.Range("N64").Value = -169162.67
$ws.Range("H65").Value = 12680.704
$ws.Range("J65").Value = 5995
$ws.Range("L65").Value = 29975
$ws.Range("N65").Value = -36215
$ws.Range("H67").Value = 65125
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 168666.67
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 168666.67
$ws.Range("M67").Value = -2142
$ws.Range("N67").Value = -170382.67
$ws.Range("H76").Value = 32261334
$ws.Range("I76").Value = 38464876
$ws.Range("J76").Value = 2920
$ws.Range("K76").Value = 38464876
$ws.Range("L76").Value = 2920
$ws.Range("M76").Value = -38464561
$ws.Range("N76").Value = -3550
$ws.Range("H79").Value = 32261334
$ws.Range("I79").Value = 38464876
$ws.Range("J79").Value = 2920
$ws.Range("K79").Value = 38464876
$ws.Range("L79").Value = 2920
$ws.Range("M79").Value = -38463784
$ws.Range("N79").Value = -5104
$ws.Range("H113").Value = 1627.8572
$ws.Range("I113").Value = 1627.8572
$ws.Range("K113").Value = 1627.8572
$ws.Range("M113").Value = 1626.1428

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3000
$ws.Range("I63").Value = 3000
$ws.Range("K63").Value = 3000
$ws.Range("M63").Value = -2314
$ws.Range("H66").Value = 3000
$ws.Range("I66").Value = 3000
$ws.Range("K66").Value = 15000
$ws.Range("M66").Value = -11568

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 32897.25
$ws.Range("J81").Value = 32897.25
$ws.Range("L81").Value = 32897.25
$ws.Range("N81").Value = -35019.25
$ws.Range("H84").Value = 32897.25
$ws.Range("J84").Value = 32897.25
$ws.Range("L84").Value = 98691.75
$ws.Range("N84").Value = -109299.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32199.697
$ws.Range("I31").Value = 39867.37
$ws.Range("J31").Value = 19260.5
$ws.Range("K31").Value = 39867.37
$ws.Range("L31").Value = 19260.5
$ws.Range("M31").Value = -39572.37
$ws.Range("N31").Value = -19850.5
$ws.Range("H34").Value = 32199.697
$ws.Range("I34").Value = 39867.37
$ws.Range("J34").Value = 19260.5
$ws.Range("K34").Value = 39867.37
$ws.Range("L34").Value = 19260.5
$ws.Range("M34").Value = -39665.37
$ws.Range("N34").Value = -19664.5
$ws.Range("H62").Value = 3153.3333
$ws.Range("J62").Value = 6500
$ws.Range("L62").Value = 6500
$ws.Range("N62").Value = -7748
$ws.Range("H65").Value = 3153.3333
$ws.Range("J65").Value = 6500
$ws.Range("L65").Value = 32500
$ws.Range("N65").Value = -38740

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 265.46155
$ws.Range("I40").Value = 265.46155
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1061.8462
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -992.8462
$ws.Range("N40").ClearContents()
$ws.Range("H68").Value = 111511.22
$ws.Range("I68").Value = 250250.25
$ws.Range("J68").Value = 520
$ws.Range("K68").Value = 750750.75
$ws.Range("L68").Value = 1560
$ws.Range("M68").Value = -749939.75
$ws.Range("N68").Value = -3182
$ws.Range("H71").Value = 111511.22
$ws.Range("I71").Value = 250250.25
$ws.Range("J71").Value = 520
$ws.Range("K71").Value = 2252252.25
$ws.Range("L71").Value = 4680
$ws.Range("M71").Value = -2248196.25
$ws.Range("N71").Value = -12792
$ws.Range("H113").Value = 514.13794
$ws.Range("I113").Value = 461.53845
$ws.Range("J113").Value = 556.875
$ws.Range("K113").Value = 1384.61535
$ws.Range("L113").Value = 1670.625
$ws.Range("M113").Value = 785.38465
$ws.Range("N113").Value = -6010.625
$ws.Range("H131").Value = 162148.55
$ws.Range("J131").Value = 173299.66
$ws.Range("L131").Value = 519898.98
$ws.Range("N131").Value = -529978.98

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H40").Value = 27614.4
$ws.Range("J40").Value = 27614.4
$ws.Range("L40").Value = 27614.4
$ws.Range("N40").Value = -27916.4
$ws.Range("H43").Value = 7285.1
$ws.Range("I43").Value = 1140.7778
$ws.Range("J43").Value = 12312.272
$ws.Range("K43").Value = 1140.7778
$ws.Range("L43").Value = 12312.272
$ws.Range("M43").Value = -989.7778000000001
$ws.Range("N43").Value = -12614.272
$ws.Range("H46").Value = 15015.333
$ws.Range("I46").Value = 5000
$ws.Range("J46").Value = 20023
$ws.Range("K46").Value = 5000
$ws.Range("L46").Value = 20023
$ws.Range("M46").Value = -4844
$ws.Range("N46").Value = -20335
$ws.Range("H57").Value = 6916.6665
$ws.Range("H70").Value = 4549396
$ws.Range("I70").Value = 6253825
$ws.Range("J70").Value = 4251.1665
$ws.Range("K70").Value = 6253825
$ws.Range("L70").Value = 4251.1665
$ws.Range("M70").Value = -6253555
$ws.Range("N70").Value = -4791.1665
$ws.Range("H73").Value = 4549396
$ws.Range("I73").Value = 6253825
$ws.Range("J73").Value = 4251.1665
$ws.Range("K73").Value = 6253825
$ws.Range("L73").Value = 4251.1665
$ws.Range("M73").Value = -6252889
$ws.Range("N73").Value = -6123.1665
$ws.Range("H80").Value = 11590.728
$ws.Range("J80").Value = 34666.668
$ws.Range("L80").Value = 34666.668
$ws.Range("N80").Value = -36662.668
$ws.Range("H83").Value = 11590.728
$ws.Range("J83").Value = 34666.668
$ws.Range("L83").Value = 173333.34
$ws.Range("N83").Value = -183317.34

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 88313
$ws.Range("J108").Value = 88313
$ws.Range("L108").Value = 88313
$ws.Range("N108").Value = -95993
